$d = $word.ActiveDocument
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" w:rsidR="00AB34F7" w:rsidRDefault="001176B0" w:rsidP="00AB34F7"><w:pPr><w:pStyle w:val="MTDisplayEquation"/></w:pPr><w:r w:rsidRPr="001176B0"><w:rPr><w:position w:val="-50"/></w:rPr><w:object w:dxaOrig="1939" w:dyaOrig="1160"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1025" type="#_x0000_t75" style="width:97.1pt;height:58pt" o:ole=""><v:imagedata r:id="rId4" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Equation.DSMT4" ShapeID="_x0000_i1025" DrawAspect="Content" ObjectID="_1707511884" r:id="rId5"/></w:object></w:r></w:p><w:p w:rsidR="00AB34F7" w:rsidRDefault="00AB34F7" w:rsidP="00AB34F7"><w:r><w:t>\[\begin{align}</w:t></w:r></w:p><w:p w:rsidR="00AB34F7" w:rsidRDefault="00AB34F7" w:rsidP="00AB34F7"><w:r><w:t xml:space="preserve">  &amp; b_{</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>0}^</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">{2}=b_{0}^{1}+(b_{1}^{1}-b_{0}^{1})*t \\ </w:t></w:r></w:p><w:p w:rsidR="00AB34F7" w:rsidRDefault="00AB34F7" w:rsidP="00AB34F7"><w:r><w:t xml:space="preserve"> &amp; ={{b}_{0</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>}}+(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">{{b}_{1}}-{{b}_{0}})*t+({{b}_{1}}-{{b}_{0}})*t+({{b}_{2}}-2{{b}_{1}}+{{b}_{0}})*{{t}^{2}} \\ </w:t></w:r></w:p><w:p w:rsidR="00AB34F7" w:rsidRDefault="00AB34F7" w:rsidP="00AB34F7"><w:r><w:t xml:space="preserve"> &amp; ={{t}^{2</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>}}{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">{b}_{0}}-2{{b}_{0}}t+{{b}_{0}}+2t(1-t){{b}_{1}}+{{t}^{2}}{{b}_{2}} \\ </w:t></w:r></w:p><w:p w:rsidR="00AB34F7" w:rsidRDefault="00AB34F7" w:rsidP="00AB34F7"><w:r><w:t xml:space="preserve"> &amp; ={{(1-t</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>)}^</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">{2}}{{b}_{0}}+2t(1-t){{b}_{1}}+{{t}^{2}}{{b}_{2}} \\ </w:t></w:r></w:p><w:p w:rsidR="00AB34F7" w:rsidRPr="00AB34F7" w:rsidRDefault="00AB34F7" w:rsidP="00AB34F7"><w:r><w:t>\end{align}\]</w:t></w:r></w:p><w:p w:rsidR="00464CA0" w:rsidRDefault="00AB34F7" w:rsidP="0046225C"><w:pPr><w:pStyle w:val="MTDisplayEquation"/></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="00AB34F7"><w:rPr><w:position w:val="-72"/></w:rPr><w:object w:dxaOrig="4700" w:dyaOrig="1560"><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:234.9pt;height:77.85pt" o:ole=""><v:imagedata r:id="rId6" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Equation.DSMT4" ShapeID="_x0000_i1026" DrawAspect="Content" ObjectID="_1707511885" r:id="rId7"/></w:object></w:r></w:p><w:p w:rsidR="0046225C" w:rsidRDefault="0046225C" w:rsidP="0046225C"><w:pPr><w:pStyle w:val="MTDisplayEquation"/></w:pPr><w:r w:rsidRPr="0046225C"><w:rPr><w:position w:val="-12"/></w:rPr><w:object w:dxaOrig="3019" w:dyaOrig="380"><v:shape id="_x0000_i1027" type="#_x0000_t75" style="width:150.85pt;height:18.9pt" o:ole=""><v:imagedata r:id="rId8" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Equation.DSMT4" ShapeID="_x0000_i1027" DrawAspect="Content" ObjectID="_1707511886" r:id="rId9"/></w:object></w:r></w:p><w:p w:rsidR="0046225C" w:rsidRDefault="00352569"><w:r w:rsidRPr="00352569"><w:rPr><w:position w:val="-30"/></w:rPr><w:object w:dxaOrig="1500" w:dyaOrig="720"><v:shape id="_x0000_i1028" type="#_x0000_t75" style="width:74.95pt;height:36.15pt" o:ole=""><v:imagedata r:id="rId10" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Equation.DSMT4" ShapeID="_x0000_i1028" DrawAspect="Content" ObjectID="_1707511887" r:id="rId11"/></w:object></w:r></w:p><w:p w:rsidR="003A5696" w:rsidRDefault="003A5696"><w:r w:rsidRPr="003A5696"><w:rPr><w:position w:val="-28"/></w:rPr><w:object w:dxaOrig="1719" w:dyaOrig="680"><v:shape id="_x0000_i1029" type="#_x0000_t75" style="width:86pt;height:33.9pt" o:ole=""><v:imagedata r:id="rId12" o:title=""/></v:shape><o:OLEObject Type="Embed" ProgID="Equation.DSMT4" ShapeID="_x0000_i1029" DrawAspect="Content" ObjectID="_1707511888" r:id="rId13"/></w:object></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:t>\[{b}'(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>0)=</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">3({{b}_{1}}-{{b}_{0}})\] </w:t></w:r></w:p><w:p><w:r><w:t>\[{b}'(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1)=</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>3({{b}_{3}}-{{b}_{2}})\]</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p>
'@
$null = $d.Content.InsertXML($xml)
